$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename headers (shared-string renames surfaced as header text changes)
$ws.Range("K1").Value = "seas_id"
$ws.Range("L1").Value = "player_id_x"
$ws.Range("N1").Value = "season_ending_year_y"
$ws.Range("O1").Value = "player_id_y"

# New player_id_y values to place in column O (rows 2-66)
$playerIdY = @{
  2 = 5062
  3 = 3962
  4 = 4528
  5 = 3187
  6 = 2223
  7 = 3372
  8 = 354
  9 = 3412
  10 = 2978
  11 = 196
  12 = 3631
  13 = 1095
  14 = 3182
  15 = 454
  16 = 5021
  17 = 1344
  18 = 3089
  19 = 614
  20 = 906
  21 = 1667
  22 = 3260
  23 = 153
  24 = 3989
  25 = 3709
  26 = 1660
  27 = 4687
  28 = 5066
  29 = 4816
  30 = 132
  31 = 1103
  32 = 1979
  33 = 2419
  34 = 917
  35 = 4575
  36 = 3226
  37 = 1337
  38 = 1239
  39 = 3758
  40 = 3485
  41 = 942
  42 = 3981
  43 = 3643
  44 = 4134
  45 = 4775
  46 = 683
  47 = 1176
  48 = 3209
  49 = 4060
  50 = 5118
  51 = 46
  52 = 145
  53 = 2330
  54 = 1706
  55 = 512
  56 = 4606
  57 = 1198
  58 = 1881
  59 = 2974
  60 = 5162
  61 = 1575
  62 = 1194
  63 = 4264
  64 = 5219
  65 = 2546
  66 = 1644
}

# Column N must hold the season_ending_year_y text (what used to be in column O);
# force text formatting before writing so strings like "2024" don't get
# silently re-typed as numbers.
$ws.Range("N2:N66").NumberFormat = "@"

for ($row = 2; $row -le 66; $row++) {
    $oldOValue = $ws.Range("O$row").Value2
    $ws.Range("N$row").Value = "$oldOValue"
    $ws.Range("O$row").Value = $playerIdY[$row]
}

# Restore column N to the workbook's normal (unstyled) cell style, now that
# the text has been committed, so no stray number-format style lingers on
# the cells themselves.
$ws.Range("N2:N66").Style = "Normal"
